$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: Model Histories data table -----------------------------------
# Implemented Bayesian Optimization for XGBoost; refreshed the tracked model
# runs (new Sklearn-Encoded / AWS-TFIDF / BO XGBOOST timings, plus the new
# "BO XGBOOST w/o Tags" and "BO XGBOOST Sold" experiments), and dropped the
# now-unused Sklearn / AWS Sagemaker rows.

# Row 3: Sklearn - Encoded
$ws1.Range("A3").Value = 1260988713616.6899
$ws1.Range("B3").Value = "Sklearn - Encoded"
$ws1.Range("C3").Value = 44262

# Row 4: AWS - TFIDF
$ws1.Range("A4").Value = 314380091392
$ws1.Range("B4").Value = "AWS - TFIDF"
$ws1.Range("C4").Value = 44255

# Row 5: BO XGBOOST
$ws1.Range("A5").Value = 1040071490897.22
$ws1.Range("B5").Value = "BO XGBOOST"
$ws1.Range("C5").Value = 44261

# Row 6: BO XGBOOST w/o Tags
$ws1.Range("A6").Value = 1125427766274.45
$ws1.Range("B6").Value = "BO XGBOOST w/o Tags"
$ws1.Range("C6").Value = 44262

# Row 7: BO XGBOOST Sold
$ws1.Range("A7").Value = 569715637462
$ws1.Range("B7").Value = "BO XGBOOST Sold"
$ws1.Range("C7").Value = 44264

# Helper formulas now only need to look as far as row 5, and the best-vs
# comparisons are re-pointed at the refreshed rows.
$ws1.Range("F3").Formula = "=VLOOKUP(MIN(A:A), A1:B5,2,FALSE)"
$ws1.Range("F5").Formula = "=1-MIN(A:A)/A3"
$ws1.Range("F6").Formula = "=1-A4/A3"

# Drop the old trailing "BO XGBOOST" run (now superseded) and leave a fresh
# blank, pre-formatted row ready for the next entry.
$ws1.Rows.Item(8).Delete()
$ws1.Range("A9").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"

# --- Sheet2: just move the saved selection ---------------------------------
$ws2.Range("F27").Select()

# Restore Sheet1 as the active/selected tab (select its range last so it
# stays the workbook's active sheet, matching the saved view state).
$ws1.Range("E13").Select()
